# Update map_info_dataframe worksheet: rename "Nation" column (country names)
# to a "City" column (specific city names), and adjust the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("A1").Value = "City"
$ws.Range("B1").Value = "color"
$ws.Range("C1").Value = "size"

# Row 2 - KERI (domestic)
$ws.Range("A2").Value = "Changon, Korea"

# Row 3 - KIMS (domestic)
$ws.Range("A3").Value = "Changon, Korea"
$ws.Range("B3").Value = "domestic"

# Row 4 - KAIST (domestic)
$ws.Range("A4").Value = "Daejeon, Korea"
$ws.Range("B4").Value = "domestic"

# Row 5 - DLR (international)
$ws.Range("A5").Value = "Koln, Germany"
$ws.Range("B5").Value = "international"

# Row 6 - University of Warwick (international)
$ws.Range("A6").Value = "Coventry, UK"
$ws.Range("B6").Value = "international"

# Row 7 - Univ of Milano bicocca (international)
$ws.Range("A7").Value = "Milan, Italy"
$ws.Range("B7").Value = "international"

# Row 8 - AIST (international)
$ws.Range("A8").Value = "Tsukuba, Japan"
$ws.Range("B8").Value = "international"

# Row 9 - NIMS (international)
$ws.Range("A9").Value = "Tsukuba, Japan"
$ws.Range("B9").Value = "international"

# Update the active cell selection to match the saved view state
$ws.Range("A9").Select()
